# Apply updated "dSF" (column F) values on Sheet1, re-pulled from source data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

$updates = @{
    2  = -4
    3  = 1
    5  = -1
    7  = -1
    8  = -1
    9  = -2
    10 = -3
    11 = -1
    12 = -2
    13 = 2
    14 = 6
    15 = -4
    16 = -3
    17 = 3
    18 = -2
    19 = -2
    20 = 0
    21 = -1
    22 = -3
    23 = 3
    24 = -2
    25 = -1
    26 = -2
    29 = 1
    30 = 0
    33 = -5
    35 = -1
    37 = -6
    39 = -3
    41 = 3
    42 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
